# [PHOENIX-5910] changes in search trade license and approval details
#
# The "searchTradeDeatils" sheet holds a sample search value used when
# searching a trade license by application number. Update that sample
# value, then make the "searchTradeDeatils" sheet the active sheet/tab
# with cell B7 selected (moving the "active" state away from the
# "licenseClosure" sheet).

$wb = $excel.ActiveWorkbook

$searchSheet = $wb.Worksheets.Item("searchTradeDeatils")

# Update the sample application-number search value.
$searchSheet.Range("B2").Value = "01987-2017-HB "

# Make this sheet the active one, with B7 selected.
$searchSheet.Activate()
$searchSheet.Range("B7").Select() | Out-Null
